$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ D = '27.738.48'; E = '  -0.91%  ' }
    3 = @{ D = '1.884.30'; E = '  -0.19%  ' }
    4 = @{ D = '1.001'; E = '  -0.13%  ' }
    5 = @{ D = '331.53'; E = '  +1.60%  ' }
    6 = @{ D = $null; E = '  -0.09%  ' }
    7 = @{ D = '0.4755'; E = '  +3.77%  ' }
    8 = @{ D = '0.3981'; E = '  +0.91%  ' }
    9 = @{ D = '48.04'; E = '  -6.29%  ' }
    10 = @{ D = '0.08060'; E = '  -1.78%  ' }
    11 = @{ D = '1.027'; E = '  -1.02%  ' }
    12 = @{ D = '21.85'; E = '  +0.82%  ' }
    13 = @{ D = '1.892.80'; E = '  -1.38%  ' }
    14 = @{ D = '5.976'; E = $null }
    15 = @{ D = '7.200'; E = '  -1.75%  ' }
    16 = @{ D = '1.001'; E = '  -0.11%  ' }
    17 = @{ D = '87.08'; E = '  -2.64%  ' }
    18 = @{ D = '0.00001045'; E = '  -1.28%  ' }
    19 = @{ D = '0.06602'; E = '  +0.56%  ' }
    20 = @{ D = '17.26'; E = '  -1.57%  ' }
    21 = @{ D = $null; E = '  -0.24%  ' }
    22 = @{ D = '27.741.16'; E = '  -0.94%  ' }
    23 = @{ D = '5.517'; E = '  -2.37%  ' }
    24 = @{ D = $null; E = '  -0.84%  ' }
    25 = @{ D = '2.308'; E = '  -0.09%  ' }
    26 = @{ D = '2.097.89'; E = '  -2.09%  ' }
    27 = @{ D = '155.93'; E = '  +1.25%  ' }
    28 = @{ D = '20.23'; E = '  +1.64%  ' }
    29 = @{ D = '2.101'; E = '  -0.49%  ' }
    30 = @{ D = '5.586'; E = '  -1.76%  ' }
    31 = @{ D = '122.52'; E = '  -1.35%  ' }
    32 = @{ D = '0.9698'; E = '  +1.14%  ' }
    33 = @{ D = '0.09558'; E = '  +0.17%  ' }
    34 = @{ D = '1.469'; E = '  -0.29%  ' }
    35 = @{ D = $null; E = '  -0.19%  ' }
    36 = @{ D = '5.310'; E = '  -2.91%  ' }
    37 = @{ D = $null; E = '  +0.24%  ' }
    38 = @{ D = '0.02257'; E = '  -1.06%  ' }
    39 = @{ D = '1.227'; E = '  -1.74%  ' }
    40 = @{ D = '8.160'; E = '  -5.47%  ' }
    41 = @{ D = '0.6011'; E = '  -1.54%  ' }
    42 = @{ D = $null; E = '  -0.10%  ' }
    43 = @{ D = $null; E = '  +0.32%  ' }
    44 = @{ D = '10.32'; E = '  -4.07%  ' }
    45 = @{ D = $null; E = '  -4.60%  ' }
    46 = @{ D = '0.5703'; E = '  -1.91%  ' }
    47 = @{ D = $null; E = '  -4.13%  ' }
    48 = @{ D = '3.410'; E = '  -0.42%  ' }
    49 = @{ D = $null; E = '  -2.55%  ' }
    50 = @{ D = '0.06821'; E = '  -0.95%  ' }
    51 = @{ D = '110.84'; E = '  +0.46%  ' }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals.D) {
        $dCell = $ws.Cells.Item([int]$row, 4)
        if ($vals.D -match '^-?\d+(\.\d+)?$') {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $vals.D
    }
    if ($null -ne $vals.E) {
        $ws.Cells.Item([int]$row, 5).Value = $vals.E
    }
}
